# Apply the "add colored cmds in console and other. weather upd" edit to the TODO list workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 3: remove the highlighted formatting on B3:D3 (match the plain style
#    used elsewhere, e.g. row 4), update task text, and clear the checkmark
#    that used to live in H3.
# ---------------------------------------------------------------------------
$ws.Range("E4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("B3").Value = "Вывести отдельную команду для реакции на грубое поведение"
$ws.Range("E3").Value = "расширить озвучивание списка дел"
$ws.Range("H3").ClearContents()

# ---------------------------------------------------------------------------
# 2. Row 4: restyle B4 like E4 (plain style instead of the bold one) and
#    update the task text for both halves of the row.
# ---------------------------------------------------------------------------
$ws.Range("E4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Доработать фильтр плохих слов"
$ws.Range("E4").Value = "почистить пылесос"

# ---------------------------------------------------------------------------
# 3. Row 5: restyle B5 like E5, then clear both halves of the row.
# ---------------------------------------------------------------------------
$ws.Range("E5").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").ClearContents()
$ws.Range("E5").ClearContents()

# ---------------------------------------------------------------------------
# 4. Rows 6-13: update task text (no style changes except B11).
# ---------------------------------------------------------------------------
$ws.Range("E6").ClearContents()
$ws.Range("B6").Value = "Сделать реакцию на восторг"
$ws.Range("B7").Value = "Сделать более медленное cэмплирование"
$ws.Range("B8").Value = "модернизировать анекдоты"
$ws.Range("B9").Value = "добавить возможность открывать ютуб"
$ws.Range("B10").Value = "поискать переключение раскладки клавиатуры на питоне и расставлять знаки препинания в зависимости от раскладки"

$ws.Range("E11").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Добавить возможность поставить пробел"

$ws.Range("B12").Value = "Улучшить фильтр для печати с клавиатуры"
$ws.Range("B13").Value = "Добавить возможность сохранять файл"

# ---------------------------------------------------------------------------
# 5. Append new rows 31-82 following the existing pattern established by
#    row 30 (same per-column styles, row height and B:D / E:G merges).
#    Rows 31-33 continue the "NN." numbering as text, rows 34-82 continue
#    as plain numbers (matching how the source sheet was authored).
# ---------------------------------------------------------------------------
for ($r = 31; $r -le 82; $r++) {
    $prev = $r - 1
    $ws.Range("B$($prev):D$($prev)").Merge() | Out-Null
    $ws.Range("E$($prev):G$($prev)").Merge() | Out-Null
}

for ($r = 31; $r -le 82; $r++) {
    $ws.Range("B$($r):D$($r)").Merge() | Out-Null
    $ws.Range("E$($r):G$($r)").Merge() | Out-Null

    $ws.Range("A30:I30").Copy()
    $ws.Range("A$($r):I$($r)").PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = $ws.Rows.Item(30).RowHeight

    if ($r -le 33) {
        $ws.Range("A$($r)").Value = "$($r - 2)."
    } else {
        $ws.Range("A$($r)").Value = $r - 2
    }
}

# ---------------------------------------------------------------------------
# 6. Restore the view selection to match the post-edit workbook.
# ---------------------------------------------------------------------------
$ws.Range("B13:D13").Select()
